# Natmi following Dr Hou advice
# Update LR-pair metrics (columns E:T, rows 2:13) on Sheet1 of the
# Dnajb11-Prtg workbook to the recomputed values (ligand/receptor cell
# counts and derived specificity/weight metrics changed after the
# pipeline re-run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.35490933333333
$ws.Range("H2").Value = 31.064728
$ws.Range("I2").Value = 0.2600176147259196
$ws.Range("J2").Value = 0.2600176147259196
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4363026666666667
$ws.Range("N2").Value = 1.308908
$ws.Range("O2").Value = 0.2059614370289619
$ws.Range("P2").Value = 0.2059614370289619
$ws.Range("Q2").Value = 4.517874555224889
$ws.Range("R2").Value = 40.660870997024
$ws.Range("S2").Value = 0.05355360158179337
$ws.Range("T2").Value = 0.05355360158179337

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.35490933333333
$ws.Range("H3").Value = 31.064728
$ws.Range("I3").Value = 0.2600176147259196
$ws.Range("J3").Value = 0.2600176147259196
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.357753333333333
$ws.Range("N3").Value = 4.073259999999999
$ws.Range("O3").Value = 0.6409422839440123
$ws.Range("P3").Value = 0.6409422839440123
$ws.Range("Q3").Value = 14.05941266369778
$ws.Range("R3").Value = 126.53471397328
$ws.Range("S3").Value = 0.1666562838481052
$ws.Range("T3").Value = 0.1666562838481052

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.35490933333333
$ws.Range("H4").Value = 31.064728
$ws.Range("I4").Value = 0.2600176147259196
$ws.Range("J4").Value = 0.2600176147259196
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3243146666666667
$ws.Range("N4").Value = 0.972944
$ws.Range("O4").Value = 0.1530962790270258
$ws.Range("P4").Value = 0.1530962790270258
$ws.Range("Q4").Value = 3.358248968803556
$ws.Range("R4").Value = 30.224240719232
$ws.Range("S4").Value = 0.03980772929602109
$ws.Range("T4").Value = 0.03980772929602109

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 13.040437
$ws.Range("H5").Value = 39.121311
$ws.Range("I5").Value = 0.3274527293839778
$ws.Range("J5").Value = 0.3274527293839778
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4363026666666667
$ws.Range("N5").Value = 1.308908
$ws.Range("O5").Value = 0.2059614370289619
$ws.Range("P5").Value = 0.2059614370289619
$ws.Range("Q5").Value = 5.689577437598667
$ws.Range("R5").Value = 51.206196938388
$ws.Range("S5").Value = 0.06744263470297986
$ws.Range("T5").Value = 0.06744263470297986

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 13.040437
$ws.Range("H6").Value = 39.121311
$ws.Range("I6").Value = 0.3274527293839778
$ws.Range("J6").Value = 0.3274527293839778
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.357753333333333
$ws.Range("N6").Value = 4.073259999999999
$ws.Range("O6").Value = 0.6409422839440123
$ws.Range("P6").Value = 0.6409422839440123
$ws.Range("Q6").Value = 17.70569680487333
$ws.Range("R6").Value = 159.35127124386
$ws.Range("S6").Value = 0.2098783002550673
$ws.Range("T6").Value = 0.2098783002550673

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 13.040437
$ws.Range("H7").Value = 39.121311
$ws.Range("I7").Value = 0.3274527293839778
$ws.Range("J7").Value = 0.3274527293839778
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3243146666666667
$ws.Range("N7").Value = 0.972944
$ws.Range("O7").Value = 0.1530962790270258
$ws.Range("P7").Value = 0.1530962790270258
$ws.Range("Q7").Value = 4.229204978842667
$ws.Range("R7").Value = 38.062844809584
$ws.Range("S7").Value = 0.05013179442593065
$ws.Range("T7").Value = 0.05013179442593065

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.797721
$ws.Range("H8").Value = 32.393163
$ws.Range("I8").Value = 0.2711368654728898
$ws.Range("J8").Value = 0.2711368654728898
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4363026666666667
$ws.Range("N8").Value = 1.308908
$ws.Range("O8").Value = 0.2059614370289619
$ws.Range("P8").Value = 0.2059614370289619
$ws.Range("Q8").Value = 4.711074466222668
$ws.Range("R8").Value = 42.399670196004
$ws.Range("S8").Value = 0.05584373844432472
$ws.Range("T8").Value = 0.05584373844432472

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.797721
$ws.Range("H9").Value = 32.393163
$ws.Range("I9").Value = 0.2711368654728898
$ws.Range("J9").Value = 0.2711368654728898
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.357753333333333
$ws.Range("N9").Value = 4.073259999999999
$ws.Range("O9").Value = 0.6409422839440123
$ws.Range("P9").Value = 0.6409422839440123
$ws.Range("Q9").Value = 14.66064168015333
$ws.Range("R9").Value = 131.94577512138
$ws.Range("S9").Value = 0.1737830818176144
$ws.Range("T9").Value = 0.1737830818176144

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.797721
$ws.Range("H10").Value = 32.393163
$ws.Range("I10").Value = 0.2711368654728898
$ws.Range("J10").Value = 0.2711368654728898
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3243146666666667
$ws.Range("N10").Value = 0.972944
$ws.Range("O10").Value = 0.1530962790270258
$ws.Range("P10").Value = 0.1530962790270258
$ws.Range("Q10").Value = 3.501859286874667
$ws.Range("R10").Value = 31.516733581872
$ws.Range("S10").Value = 0.04151004521095071
$ws.Range("T10").Value = 0.04151004521095071

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5.630808999999999
$ws.Range("H11").Value = 16.892427
$ws.Range("I11").Value = 0.1413927904172128
$ws.Range("J11").Value = 0.1413927904172128
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.4363026666666667
$ws.Range("N11").Value = 1.308908
$ws.Range("O11").Value = 0.2059614370289619
$ws.Range("P11").Value = 0.2059614370289619
$ws.Range("Q11").Value = 2.456736982190666
$ws.Range("R11").Value = 22.110632839716
$ws.Range("S11").Value = 0.02912146229986397
$ws.Range("T11").Value = 0.02912146229986397

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 5.630808999999999
$ws.Range("H12").Value = 16.892427
$ws.Range("I12").Value = 0.1413927904172128
$ws.Range("J12").Value = 0.1413927904172128
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.357753333333333
$ws.Range("N12").Value = 4.073259999999999
$ws.Range("O12").Value = 0.6409422839440123
$ws.Range("P12").Value = 0.6409422839440123
$ws.Range("Q12").Value = 7.645249689113331
$ws.Range("R12").Value = 68.80724720201998
$ws.Range("S12").Value = 0.0906246180232254
$ws.Range("T12").Value = 0.0906246180232254

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 5.630808999999999
$ws.Range("H13").Value = 16.892427
$ws.Range("I13").Value = 0.1413927904172128
$ws.Range("J13").Value = 0.1413927904172128
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.3243146666666667
$ws.Range("N13").Value = 0.972944
$ws.Range("O13").Value = 0.1530962790270258
$ws.Range("P13").Value = 0.1530962790270258
$ws.Range("Q13").Value = 1.826153943898667
$ws.Range("R13").Value = 16.435385495088
$ws.Range("S13").Value = 0.02164671009412339
$ws.Range("T13").Value = 0.02164671009412339
